$wb = $excel.ActiveWorkbook

# Update status text "Ready for handoff" -> "In Translation" on all sheets/cells
# where it appears: Overview!E2, Overview!F2 (zh-cn/de-de status on the Overview
# sheet), and the per-locale "Status" column (column C) on zh-cn / de-de.
# Shrinking the status text re-triggers the report generator's column-width
# fit for the columns that display it, so those columns narrow as well.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
